# ThaiModelInput.xlsx - "Added 2050 with pink upgrades"
#
# Changes applied:
#  - lines sheet: max_i_ka (column G, rows 2-21) raised from 40 to 200
#  - selections updated on each sheet to reflect where the author was
#    last working, with the "lines" sheet ending up as the active tab
#    (selection G2:G21, the cells that were just edited)

$wb = $excel.ActiveWorkbook

$wsBuses = $wb.Worksheets.Item(1)   # buses
$wsLines = $wb.Worksheets.Item(2)   # lines
$wsLoads = $wb.Worksheets.Item(3)   # loads
$wsGen   = $wb.Worksheets.Item(4)   # gen

# Bump the max_i_ka values for every line row (2-21) from 40 -> 200
$wsLines.Range("G2:G21").Value = 200

# Update the remembered selection on sheets that aren't ending up active
$wsBuses.Range("I3").Select()
$wsLoads.Range("C19").Select()
$wsGen.Range("D11").Select()

# Finish on the "lines" sheet with the edited range selected; this also
# makes it the active tab, matching the saved workbook view state.
$wsLines.Activate()
$wsLines.Range("G2:G21").Select()
